# MHD2-138: Clinical context update - APL (12Sep2024)
# The References bibliography paragraph (ADDIN EN.REFLIST field) switches
# from the EndNoteBibliography style (with an explicit both-justify) to the
# document's own CLIN4 style, and the EndNote-generated "noProof" /
# complex-script-font run properties that EndNote stamped onto every run
# are stripped out (CLIN4 already carries equivalent formatting via the
# style definition).

$d = $word.ActiveDocument

# Locate the References bibliography paragraph: the one whose style is
# "EndNote Bibliography" (w:pStyle val="EndNoteBibliography").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Style.NameLocal -eq "EndNote Bibliography") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the EndNote Bibliography paragraph"
}

$range = $target.Range

# Pull this paragraph's slice of OOXML (WordOpenXML always returns/accepts
# a full single-part package wrapper, so edit the fragment in place inside
# that wrapper rather than trying to round-trip a bare <w:p>).
$xml = $range.WordOpenXML

$startTag = '<w:p w14:paraId="04E6F810"'
$startIdx = $xml.IndexOf($startTag)
if ($startIdx -lt 0) {
    throw "Could not locate target paragraph start in WordOpenXML"
}
$endTag = "</w:p>"
$endIdx = $xml.IndexOf($endTag, $startIdx) + $endTag.Length

$paraXml = $xml.Substring($startIdx, $endIdx - $startIdx)
$original = $paraXml

# 1. Paragraph formatting: drop the EndNoteBibliography style + explicit
#    both-justification in favour of the CLIN4 paragraph style (CLIN4
#    already bakes in jc="both").
$paraXml = $paraXml.Replace(
    '<w:pPr><w:pStyle w:val="EndNoteBibliography"/><w:jc w:val="both"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="CLIN4"/></w:pPr>'
)

# 2. Strip the complex-script font override that only wrapped the two
#    ADDIN EN.REFLIST fldChar runs.
$paraXml = $paraXml.Replace('<w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr>', '')

# 3. Remove the EndNote-stamped w:noProof on every bibliography run; once
#    removed, collapse any run-properties element left empty.
$paraXml = $paraXml.Replace('<w:noProof/>', '')
$paraXml = $paraXml.Replace('<w:rPr></w:rPr>', '')

if ($paraXml -eq $original) {
    throw "No changes were made to the target paragraph"
}

$newXml = $xml.Substring(0, $startIdx) + $paraXml + $xml.Substring($endIdx)

$range.WordOpenXML = $newXml

Write-Host "Updated bibliography paragraph style and stripped EndNote run formatting."
